# Generate Report for Handback
#
# This applies the "handback" results to the localization-status workbook:
#   - Overview sheet: status text updates from "Ready for handoff" to
#     "Handed back: in sync with en-US" (columns E/F), and those two
#     columns are widened to fit the longer text.
#   - zh-cn / de-de sheets: the "Latest Target File" / "Latest Handback
#     File" / "Latest Handback DateTime" columns (I/J/K) are populated
#     now that handback has happened, including turning the target-file
#     cell into a hyperlink back to the source markdown doc. Those
#     columns are widened as well.

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/602696415cc9f0f49c9236a386b14856d8f248cd/e2e/"
$targetMd = "ec347a14-2d49-483f-b70a-b82f66d3c54f.md"
$targetUrl = $baseUrl + $targetMd

# --- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# Widen columns E and F so the longer status text is readable.
$overview.Columns.Item(5).ColumnWidth = 29.1667
$overview.Columns.Item(6).ColumnWidth = 29.1667

# --- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

# Latest Target File (I) becomes a hyperlink to the source doc.
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $targetUrl, $null, $null, $targetMd)
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $targetUrl, $null, $null, $targetMd)

# Latest Handback File (J) + Latest Handback DateTime (K).
$zhcn.Range("J2").Value = "ec347a14-2d49-483f-b70a-b82f66d3c54f.5b45b90031a49edabfb6b20291b1375b956b5d29.zh-cn.xlf"
$zhcn.Range("J3").Value = "ec347a14-2d49-483f-b70a-b82f66d3c54f.5b45b90031a49edabfb6b20291b1375b956b5d29.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-29 13:03:49"
$zhcn.Range("K3").Value = "2016-08-29 13:03:49"

# Widen columns C, I, J to fit the new content.
$zhcn.Columns.Item(3).ColumnWidth = 29.1667
$zhcn.Columns.Item(9).ColumnWidth = 39.1667
$zhcn.Columns.Item(10).ColumnWidth = 39.1667

# --- de-de sheet --------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Hyperlinks.Add($dede.Range("I2"), $targetUrl, $null, $null, $targetMd)
$dede.Hyperlinks.Add($dede.Range("I3"), $targetUrl, $null, $null, $targetMd)

$dede.Range("J2").Value = "ec347a14-2d49-483f-b70a-b82f66d3c54f.5b45b90031a49edabfb6b20291b1375b956b5d29.de-de.xlf"
$dede.Range("J3").Value = "ec347a14-2d49-483f-b70a-b82f66d3c54f.5b45b90031a49edabfb6b20291b1375b956b5d29.de-de.xlf"
$dede.Range("K2").Value = "2016-08-29 13:03:56"
$dede.Range("K3").Value = "2016-08-29 13:03:56"

$dede.Columns.Item(3).ColumnWidth = 29.1667
$dede.Columns.Item(9).ColumnWidth = 39.1667
$dede.Columns.Item(10).ColumnWidth = 39.1667
